$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @{Row=2;  B=110000; D=0.3646489200594908; E=6.143895547706342;  F=0.8832;  H=6.423239529546289},
    @{Row=3;  B=110000; D=0.7034468968829464; E=5.790383126219522;  F=1.6785;  H=6.423239529546289},
    @{Row=4;  B=110000; D=0.9685449757105207; E=5.585202231743557;  F=2.1185;  H=6.423239529546289},
    @{Row=5;  B=110000; D=1.13743123501128;   E=5.468366089922898;  F=2.6534;  H=6.423239529546289},
    @{Row=6;  B=110000; D=1.201221854365843;  E=5.426089715626122;  F=2.9961;  H=6.423239529546289},
    @{Row=7;  B=110000; D=1.233482886920558;  E=5.436908337885208;  F=3.3062;  H=6.423239529546289},
    @{Row=8;  B=110000; D=1.268277853131055;  E=5.420993439048613;  F=3.4356;  H=6.423239529546289},
    @{Row=9;  B=110000; D=1.290522950475432;  E=5.410628427984358;  F=3.578;   H=6.423239529546289},
    @{Row=10; B=110000; D=1.316234812580474;  E=5.38875204158206;   F=3.6839;  H=6.423239529546289}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 8).Value = $item.H
}
